# This edit reorders the data rows (2-25) of the sheet: columns D (Fecha),
# K (Variedad), L (Calidad), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado), Q (Unidad de comercialización), R (Origen),
# S (Precio $/Kg) and T (Kg / unidad) get shuffled between rows according to
# the mapping below (columns A, B, C, E, F, G, H, I, J stay constant/unchanged).
#
# Mapping key = destination (new) row number, value = source (original) row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 6
    3  = 2
    4  = 3
    5  = 16
    6  = 17
    7  = 14
    8  = 15
    9  = 24
    10 = 25
    11 = 4
    12 = 5
    13 = 19
    14 = 20
    15 = 21
    16 = 7
    17 = 8
    18 = 9
    19 = 22
    20 = 23
    21 = 12
    22 = 13
    23 = 18
    24 = 10
    25 = 11
}

# Columns whose values move together with the row when shuffling.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# First capture the original values for every relevant cell so that writes
# performed while iterating do not clobber values still needed as a source.
# `.Value2` is used (rather than `.Value`) so dates come back as the raw
# numeric serial number, matching the underlying stored cell value.
$original = @{}
foreach ($r in 2..25) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Now write the shuffled values into their destination rows.
foreach ($destRow in 2..25) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
